# Adds mobility-data columns (Retail, Grocery, Parks, Transit, Workplace,
# Residential) to the Indonesia shapefile-source sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1), columns E-J -----------------------------------
$headers = @("Retail","Grocery","Parks","Transit","Workplace","Residential")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, 5 + $i).Value = $headers[$i]
}

# --- Mobility data for rows 2-35, columns E-J ---------------------------
$data = @(
    @(-2,8,-11,-28,-10,6),
    @(-34,-26,-43,-69,-37,15),
    @(-13,9,-9,-36,-15,9),
    @(-11,-6,-4,-47,-23,14),
    @(-1,8,-3,-30,-16,8),
    @(-6,1,-9,-37,-18,13),
    @(-14,-6,-16,-27,-16,9),
    @(-1,6,1,-29,-11,8),
    @(-10,-4,-14,-38,-21,14),
    @(-13,-1,-10,-34,-19,11),
    @(-7,1,2,-19,-7,7),
    @(-4,-2,21,0,-17,13),
    @(-26,-8,-71,-42,-31,14),
    @(-11,3,-12,-27,-14,8),
    @(-3,13,-8,-29,-20,9),
    @(-20,-11,-14,-33,-20,11),
    @(-12,9,-7,-29,-11,8),
    @(-21,9,-20,-9,-18,10),
    @(-27,-8,-20,-45,-22,13),
    @(-20,-11,-16,-36,-22,11),
    @(-20,-11,-10,-31,-14,9),
    @(-9,-3,-5,-41,-15,9),
    @(-20,-12,-25,-50,-15,11),
    @(-3,8,-7,-28,-13,9),
    @(-17,-7,-10,-37,-19,12),
    @(-15,-2,-7,-47,-25,12),
    @(-13,1,5,-41,-19,10),
    @(-23,-17,-32,-56,-29,13),
    @(-23,3,-2,-28,-23,13),
    @(-23,-12,-18,-36,-18,10),
    @(-13,-4,-17,-53,-15,12),
    @(-20,-7,-31,-47,-19,11),
    @(1,13,2,-29,-12,10),
    @(-7,4,-9,-31,-19,8),
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item(2 + $r, 5 + $c).Value = $row[$c]
    }
}

# --- Column sizing --------------------------------------------------
# D ("Recovered"), I ("Workplace") and J ("Residential") are the widest
# columns in the new layout, so Excel auto-fits them to the header text.
$ws.Range("D1:D35").EntireColumn.AutoFit() | Out-Null
$ws.Range("I1:J35").EntireColumn.AutoFit() | Out-Null

# --- Selection moves to H1 (top of the new Transit column) -------------
$ws.Range("H1").Select() | Out-Null

Write-Output "done"
